$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly refresh of "Hortaliza, Terminal La Palmera de La Serena - Ciboulette":
# each data row (2-31) is reassigned the Fecha / Volumen / Precio minimo /
# Precio maximo / Precio promedio ponderado / Precio por Kg values that
# originally belonged to another row in the sheet (a full permutation of
# those six columns across the 30 data rows).
$data = @{
    2  = @(44978, 1000, 1800, 2000, 1900, 633)
    3  = @(44965, 1120, 2000, 2500, 2250, 750)
    4  = @(45006, 1100, 2000, 2500, 2250, 750)
    5  = @(44848, 1000, 1500, 2000, 1750, 583)
    6  = @(45020, 1200, 2000, 2500, 2250, 750)
    7  = @(44985, 1000, 2000, 2500, 2250, 750)
    8  = @(44881, 500, 1900, 2000, 1950, 650)
    9  = @(45070, 800, 2000, 2500, 2250, 750)
    10 = @(44999, 1100, 2000, 2500, 2250, 750)
    11 = @(44971, 1000, 2000, 2500, 2250, 750)
    12 = @(44992, 1040, 2000, 2500, 2250, 750)
    13 = @(45062, 1100, 2000, 2500, 2250, 750)
    14 = @(44827, 1200, 2000, 2500, 2250, 750)
    15 = @(44911, 700, 1800, 2000, 1900, 633)
    16 = @(44970, 800, 2000, 2500, 2250, 750)
    17 = @(45035, 1100, 2000, 2500, 2250, 750)
    18 = @(45034, 1100, 2000, 2500, 2250, 750)
    19 = @(44964, 1000, 2000, 2500, 2250, 750)
    20 = @(45041, 1160, 2000, 2500, 2250, 750)
    21 = @(45028, 1000, 2000, 2500, 2250, 750)
    22 = @(45013, 1100, 2000, 2500, 2250, 750)
    23 = @(44685, 400, 1500, 2000, 1750, 583)
    24 = @(45084, 900, 2000, 2500, 2250, 750)
    25 = @(44883, 500, 1800, 2000, 1900, 633)
    26 = @(44951, 800, 2000, 2500, 2250, 750)
    27 = @(45091, 800, 2000, 2500, 2250, 750)
    28 = @(45007, 1160, 2000, 2500, 2250, 750)
    29 = @(44910, 1000, 1800, 2000, 1900, 633)
    30 = @(45077, 760, 2000, 2500, 2250, 750)
    31 = @(44953, 1000, 2000, 2500, 2250, 750)
}

for ($row = 2; $row -le 31; $row++) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]    # D: Fecha
    $ws.Cells.Item($row, 10).Value = $vals[1]   # J: Volumen
    $ws.Cells.Item($row, 11).Value = $vals[2]   # K: Precio minimo
    $ws.Cells.Item($row, 12).Value = $vals[3]   # L: Precio maximo
    $ws.Cells.Item($row, 13).Value = $vals[4]   # M: Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value = $vals[5]   # P: Precio por Kg
}
